# Weekly update: insert this week's price rows (new week, date 45239) above the
# previous week's rows (date 45173) for "Agrícola del Norte S.A. de Arica - Frutilla".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing rows 150:153 down by 4 rows (they become rows 154:157)
# and open up blank rows 150:153 for the new week's data.
$ws.Rows("150:153").Insert()

# Shared, unchanging attributes for every one of these rows.
$mercado   = "Agrícola del Norte S.A. de Arica"
$region    = "Arica y Parinacota"
$tipo      = "Fruta"
$producto  = "Berries"
$categoria = "Frutilla"
$variedad  = "Sin especificar"
$unidad    = "$/bandeja 3 kilos"
$origen    = "Región de Arica y Parinacota"

$calidades = @("Especial", "Primera", "Segunda", "Tercera")

# New week's data (2023-11-09, serial 45239).
$nuevaFecha = 45239
$volumen    = @(220, 250, 250, 110)
$precioMin  = @(7000, 5000, 3000, 2000)
$precioMax  = @(8000, 6000, 4000, 3000)
$precioProm = @(7455, 5600, 3600, 2545)
$precioKg   = @(2485, 1867, 1200, 848)

for ($i = 0; $i -lt 4; $i++) {
    $r = 150 + $i
    $ws.Cells.Item($r, 1).Value = 1
    $ws.Cells.Item($r, 2).Value = $mercado
    $ws.Cells.Item($r, 3).Value = $region
    $ws.Cells.Item($r, 4).Value = $nuevaFecha
    $ws.Cells.Item($r, 5).Value = 15
    $ws.Cells.Item($r, 6).Value = $tipo
    $ws.Cells.Item($r, 7).Value = 100101
    $ws.Cells.Item($r, 8).Value = $producto
    $ws.Cells.Item($r, 9).Value = 100112025
    $ws.Cells.Item($r, 10).Value = $categoria
    $ws.Cells.Item($r, 11).Value = $variedad
    $ws.Cells.Item($r, 12).Value = $calidades[$i]
    $ws.Cells.Item($r, 13).Value = $volumen[$i]
    $ws.Cells.Item($r, 14).Value = $precioMin[$i]
    $ws.Cells.Item($r, 15).Value = $precioMax[$i]
    $ws.Cells.Item($r, 16).Value = $precioProm[$i]
    $ws.Cells.Item($r, 17).Value = $unidad
    $ws.Cells.Item($r, 18).Value = $origen
    $ws.Cells.Item($r, 19).Value = $precioKg[$i]
    $ws.Cells.Item($r, 20).Value = 3
}
